$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 57.38695133333332
$ws.Range("H2").Value = 172.160854
$ws.Range("I2").Value = 0.6542464432660997
$ws.Range("J2").Value = 0.6542464432660998
$ws.Range("M2").Value = 181.3526613333333
$ws.Range("N2").Value = 544.057984
$ws.Range("O2").Value = 0.9845849379007657
$ws.Range("P2").Value = 0.984584937900766
$ws.Range("Q2").Value = 10407.27635010648
$ws.Range("R2").Value = 93665.48715095833
$ws.Range("S2").Value = 0.6441611937149496
$ws.Range("T2").Value = 0.6441611937149498
$ws.Range("G3").Value = 57.38695133333332
$ws.Range("H3").Value = 172.160854
$ws.Range("I3").Value = 0.6542464432660997
$ws.Range("J3").Value = 0.6542464432660998
$ws.Range("O3").Value = 0.003686045149950483
$ws.Range("P3").Value = 0.003686045149950484
$ws.Range("Q3").Value = 38.96229673825332
$ws.Range("R3").Value = 350.6606706442799
$ws.Range("S3").Value = 0.002411581929073361
$ws.Range("T3").Value = 0.002411581929073362
$ws.Range("G4").Value = 57.38695133333332
$ws.Range("H4").Value = 172.160854
$ws.Range("I4").Value = 0.6542464432660997
$ws.Range("J4").Value = 0.6542464432660998
$ws.Range("M4").Value = 0.6398506666666667
$ws.Range("N4").Value = 1.919552
$ws.Range("O4").Value = 0.003473824559694892
$ws.Range("P4").Value = 0.003473824559694892
$ws.Range("Q4").Value = 36.71907906860088
$ws.Range("R4").Value = 330.4717116174079
$ws.Range("S4").Value = 0.002272737362710808
$ws.Range("T4").Value = 0.002272737362710808
$ws.Range("G5").Value = 57.38695133333332
$ws.Range("H5").Value = 172.160854
$ws.Range("I5").Value = 0.6542464432660997
$ws.Range("J5").Value = 0.6542464432660998
$ws.Range("M5").Value = 1.520540333333333
$ws.Range("N5").Value = 4.561621
$ws.Range("O5").Value = 0.008255192389588805
$ws.Range("P5").Value = 0.008255192389588807
$ws.Range("Q5").Value = 87.25917410937042
$ws.Range("R5").Value = 785.3325669843338
$ws.Range("S5").Value = 0.005400930259365851
$ws.Range("T5").Value = 0.005400930259365853
$ws.Range("I6").Value = 0.16357689713892
$ws.Range("J6").Value = 0.16357689713892
$ws.Range("M6").Value = 181.3526613333333
$ws.Range("N6").Value = 544.057984
$ws.Range("O6").Value = 0.9845849379007657
$ws.Range("P6").Value = 0.984584937900766
$ws.Range("Q6").Value = 2602.062251220025
$ws.Range("R6").Value = 23418.56026098023
$ws.Range("S6").Value = 0.1610553491115234
$ws.Range("T6").Value = 0.1610553491115235
$ws.Range("I7").Value = 0.16357689713892
$ws.Range("J7").Value = 0.16357689713892
$ws.Range("O7").Value = 0.003686045149950483
$ws.Range("P7").Value = 0.003686045149950484
$ws.Range("S7").Value = 0.000602951828342865
$ws.Range("T7").Value = 0.0006029518283428651
$ws.Range("I8").Value = 0.16357689713892
$ws.Range("J8").Value = 0.16357689713892
$ws.Range("M8").Value = 0.6398506666666667
$ws.Range("N8").Value = 1.919552
$ws.Range("O8").Value = 0.003473824559694892
$ws.Range("P8").Value = 0.003473824559694892
$ws.Range("Q8").Value = 9.180627700252446
$ws.Range("R8").Value = 82.62564930227201
$ws.Range("S8").Value = 0.0005682374426798653
$ws.Range("T8").Value = 0.0005682374426798653
$ws.Range("I9").Value = 0.16357689713892
$ws.Range("J9").Value = 0.16357689713892
$ws.Range("M9").Value = 1.520540333333333
$ws.Range("N9").Value = 4.561621
$ws.Range("O9").Value = 0.008255192389588805
$ws.Range("P9").Value = 0.008255192389588807
$ws.Range("Q9").Value = 21.81683231850622
$ws.Range("R9").Value = 196.351490866556
$ws.Range("S9").Value = 0.001350358756373763
$ws.Range("T9").Value = 0.001350358756373763
$ws.Range("G10").Value = 14.516389
$ws.Range("H10").Value = 43.549167
$ws.Range("I10").Value = 0.1654957381714162
$ws.Range("J10").Value = 0.1654957381714162
$ws.Range("M10").Value = 181.3526613333333
$ws.Range("N10").Value = 544.057984
$ws.Range("O10").Value = 0.9845849379007657
$ws.Range("P10").Value = 0.984584937900766
$ws.Range("Q10").Value = 2632.585778099925
$ws.Range("R10").Value = 23693.27200289933
$ws.Range("S10").Value = 0.1629446110903452
$ws.Range("T10").Value = 0.1629446110903452
$ws.Range("G11").Value = 14.516389
$ws.Range("H11").Value = 43.549167
$ws.Range("I11").Value = 0.1654957381714162
$ws.Range("J11").Value = 0.1654957381714162
$ws.Range("O11").Value = 0.003686045149950483
$ws.Range("P11").Value = 0.003686045149950484
$ws.Range("Q11").Value = 9.855757147659999
$ws.Range("R11").Value = 88.70181432894
$ws.Range("S11").Value = 0.0006100247630242236
$ws.Range("T11").Value = 0.0006100247630242237
$ws.Range("G12").Value = 14.516389
$ws.Range("H12").Value = 43.549167
$ws.Range("I12").Value = 0.1654957381714162
$ws.Range("J12").Value = 0.1654957381714162
$ws.Range("M12").Value = 0.6398506666666667
$ws.Range("N12").Value = 1.919552
$ws.Range("O12").Value = 0.003473824559694892
$ws.Range("P12").Value = 0.003473824559694892
$ws.Range("Q12").Value = 9.288321179242665
$ws.Range("R12").Value = 83.59489061318399
$ws.Range("S12").Value = 0.0005749031597847009
$ws.Range("T12").Value = 0.000574903159784701
$ws.Range("G13").Value = 14.516389
$ws.Range("H13").Value = 43.549167
$ws.Range("I13").Value = 0.1654957381714162
$ws.Range("J13").Value = 0.1654957381714162
$ws.Range("M13").Value = 1.520540333333333
$ws.Range("N13").Value = 4.561621
$ws.Range("O13").Value = 0.008255192389588805
$ws.Range("P13").Value = 0.008255192389588807
$ws.Range("Q13").Value = 22.07275496885633
$ws.Range("R13").Value = 198.654794719707
$ws.Range("S13").Value = 0.001366199158262056
$ws.Range("T13").Value = 0.001366199158262056
$ws.Range("G14").Value = 1.46316
$ws.Range("H14").Value = 4.389480000000001
$ws.Range("I14").Value = 0.01668092142356404
$ws.Range("J14").Value = 0.01668092142356404
$ws.Range("M14").Value = 181.3526613333333
$ws.Range("N14").Value = 544.057984
$ws.Range("O14").Value = 0.9845849379007657
$ws.Range("P14").Value = 0.984584937900766
$ws.Range("Q14").Value = 265.3479599564801
$ws.Range("R14").Value = 2388.131639608321
$ws.Range("S14").Value = 0.01642378398394735
$ws.Range("T14").Value = 0.01642378398394735
$ws.Range("G15").Value = 1.46316
$ws.Range("H15").Value = 4.389480000000001
$ws.Range("I15").Value = 0.01668092142356404
$ws.Range("J15").Value = 0.01668092142356404
$ws.Range("O15").Value = 0.003686045149950483
$ws.Range("P15").Value = 0.003686045149950484
$ws.Range("Q15").Value = 0.9933978504000002
$ws.Range("R15").Value = 8.940580653600001
$ws.Range("S15").Value = 0.00006148662951003332
$ws.Range("T15").Value = 0.00006148662951003335
$ws.Range("G16").Value = 1.46316
$ws.Range("H16").Value = 4.389480000000001
$ws.Range("I16").Value = 0.01668092142356404
$ws.Range("J16").Value = 0.01668092142356404
$ws.Range("M16").Value = 0.6398506666666667
$ws.Range("N16").Value = 1.919552
$ws.Range("O16").Value = 0.003473824559694892
$ws.Range("P16").Value = 0.003473824559694892
$ws.Range("Q16").Value = 0.9362039014400002
$ws.Range("R16").Value = 8.425835112960002
$ws.Range("S16").Value = 0.00005794659451951743
$ws.Range("T16").Value = 0.00005794659451951745
$ws.Range("G17").Value = 1.46316
$ws.Range("H17").Value = 4.389480000000001
$ws.Range("I17").Value = 0.01668092142356404
$ws.Range("J17").Value = 0.01668092142356404
$ws.Range("M17").Value = 1.520540333333333
$ws.Range("N17").Value = 4.561621
$ws.Range("O17").Value = 0.008255192389588805
$ws.Range("P17").Value = 0.008255192389588807
$ws.Range("Q17").Value = 2.22479379412
$ws.Range("R17").Value = 20.02314414708
$ws.Range("S17").Value = 0.0001377042155871347
$ws.Range("T17").Value = 0.0001377042155871348
